$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update a few odds in row 2 (existing match, odds tweaked) ---
$ws.Range("G2").Value = 3.3
$ws.Range("I2").Value = 2.45
$ws.Range("L2").Value = 3.2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("W2").Value = 8.5
$ws.Range("AC2").Value = 7
$ws.Range("AI2").Value = 11
$ws.Range("AJ2").Value = 10
$ws.Range("AK2").Value = 23
$ws.Range("AX2").Value = 15
$ws.Range("AZ2").Value = 51
$ws.Range("BA2").Value = 81

# --- Step 2: insert a new row at position 6, shifting old row 6 (and below) down to row 7 ---
$ws.Rows("6:6").Insert()

# --- Step 3: populate the newly inserted row 6 with the new match data ---
$ws.Range("A6").Value = "4ECJO79E"
# B6's date string matches the other rows' Date column exactly; copy an existing
# text-typed cell instead of assigning a literal, so Excel doesn't reinterpret
# "07/11/2024" as a date serial number (keeps it stored as plain text, like the rest).
$ws.Range("B2").Copy($ws.Range("B6"))
$ws.Range("C6").Value = "20:45"
$ws.Range("D6").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E6").Value = "Dep. Pasto"
$ws.Range("F6").Value = "America De Cali"
$ws.Range("G6").Value = 1.95
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 2.75
$ws.Range("K6").Value = 1.95
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("S6").Value = 1.57
$ws.Range("T6").Value = 2.25
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("W6").Value = 5.5
$ws.Range("X6").Value = 8
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 17
$ws.Range("AA6").Value = 19
$ws.Range("AB6").Value = 41
$ws.Range("AC6").Value = 6.5
$ws.Range("AD6").Value = 6
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 81
$ws.Range("AG6").Value = 201
$ws.Range("AH6").Value = 9.5
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 17
$ws.Range("AK6").Value = 51
$ws.Range("AL6").Value = 41
$ws.Range("AM6").Value = 51
$ws.Range("AN6").Value = 3.75
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 26
$ws.Range("AQ6").Value = 41
$ws.Range("AR6").Value = 67
$ws.Range("AS6").Value = 251
$ws.Range("AT6").Value = 2.25
$ws.Range("AU6").Value = 9.5
$ws.Range("AV6").Value = 81
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 26
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 101
$ws.Range("BA6").Value = 151
$ws.Range("BB6").Value = 351
$ws.Range("BC6").Value = 126
$ws.Range("BD6").Value = 126

# --- Step 4: update a few odds in row 7 (the shifted-down original row, odds tweaked) ---
$ws.Range("G7").Value = 2.4
$ws.Range("H7").Value = 3.1
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.38
$ws.Range("R7").Value = 1.57
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.73
$ws.Range("W7").Value = 6.5
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 451
$ws.Range("AJ7").Value = 12
$ws.Range("AL7").Value = 29
$ws.Range("AS7").Value = 251
$ws.Range("AX7").Value = 19
$ws.Range("AZ7").Value = 67
$ws.Range("BA7").Value = 101

Write-Output "edit applied"